$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Step 1: Fill in the "Fredag, 1. april" day-log entry. In the
# original document this was an empty placeholder paragraph right
# after the "Fredag, 1. april" Heading3 paragraph; we give it the
# day's notes and append two more paragraphs for the work location
# and hours, mirroring the pattern used by the other day entries.
# ------------------------------------------------------------------
$p = $d.Paragraphs.Item(59)
$r = $p.Range
$r.Text = "I dag har jeg fået kortet helt til at virke, og fremvise ruter, samt start og slut punkter. Ydermere har jeg fået deployet appen til firebase og kontrolleret at den køre tilfredsstillende på telefonen."
$r.InsertParagraphAfter()

$p2 = $d.Paragraphs.Item(60)
$r2 = $p2.Range
$r2.Text = "Arbejdssted: Hjemme, Bjerringbro."
$r2.InsertParagraphAfter()

$p3 = $d.Paragraphs.Item(61)
$r3 = $p3.Range
$r3.Text = "Arbejdstimer: 7."

# ------------------------------------------------------------------
# Step 2: The extra content above pushes the page break one heading
# earlier, so move the cached <w:lastRenderedPageBreak/> marker from
# "Søndag, 10. april" to "Fredag, 8. april".
# ------------------------------------------------------------------

# 2a. Add <w:lastRenderedPageBreak/> to the run of "Fredag, 8. april"
#     (now at paragraph index 75 after the 2 new paragraphs above).
$pFre8 = $d.Paragraphs.Item(75)
$fullFre8 = $pFre8.Range
$runFre8 = $d.Range($fullFre8.Start, $fullFre8.End - 1)
$xmlFre8 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r w:rsidRPr="00FF19BA"><w:lastRenderedPageBreak/><w:t>Fredag, 8. april</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$runFre8.InsertXML($xmlFre8)

# 2b. Remove <w:lastRenderedPageBreak/> from the run of "Søndag, 10. april"
#     (now at paragraph index 79).
$pSon10 = $d.Paragraphs.Item(79)
$fullSon10 = $pSon10.Range
$runSon10 = $d.Range($fullSon10.Start, $fullSon10.End - 1)
$xmlSon10 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r w:rsidRPr="00FF19BA"><w:t>Søndag, 10. april</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$runSon10.InsertXML($xmlSon10)

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
